# Updates the crypto-price table to the figures from the latest GitHub
# Actions refresh (prices/volumes refreshed, two coins re-ranked/swapped).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.144.07'
$ws.Range("E2").Value = '  -0.96%  '
$ws.Range("D3").Value = '1.846.18'
$ws.Range("E3").Value = '  -2.40%  '
$ws.Range("D4").Value = '''1.000'
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = '''235.98'
$ws.Range("E5").Value = '  -0.88%  '
$ws.Range("D6").Value = '''1.000'
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("D7").Value = '''0.4777'
$ws.Range("E7").Value = '  -2.50%  '
$ws.Range("D8").Value = '''0.2803'
$ws.Range("E8").Value = '  -4.61%  '
$ws.Range("D9").Value = '''0.06468'
$ws.Range("E9").Value = '  -3.49%  '
$ws.Range("D10").Value = '1.857.80'
$ws.Range("E10").Value = '  -1.90%  '
$ws.Range("D11").Value = '''0.07313'
$ws.Range("E11").Value = '  -0.50%  '
$ws.Range("D12").Value = '''16.25'
$ws.Range("E12").Value = '  -4.28%  '
$ws.Range("D13").Value = '''5.105'
$ws.Range("E13").Value = '  -0.54%  '
$ws.Range("D14").Value = '''87.08'
$ws.Range("E14").Value = '  -0.90%  '
$ws.Range("E15").Value = '  -2.89%  '
$ws.Range("D16").Value = '30.090.92'
$ws.Range("E16").Value = '  -1.09%  '
$ws.Range("E17").Value = '  +0.06%  '
$ws.Range("D18").Value = '''13.20'
$ws.Range("E18").Value = '  -1.80%  '
$ws.Range("D19").Value = '''0.000007615'
$ws.Range("E19").Value = '  -2.58%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = '''223.93'
$ws.Range("E20").Value = '  +17.82%  '
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.098.93'
$ws.Range("E21").Value = '  -1.97%  '
$ws.Range("D22").Value = '''1.001'
$ws.Range("E22").Value = '  -0.18%  '
$ws.Range("D23").Value = '''5.281'
$ws.Range("E23").Value = '  -0.48%  '
$ws.Range("D24").Value = '''6.068'
$ws.Range("E24").Value = '  -1.14%  '
$ws.Range("D25").Value = '''9.207'
$ws.Range("E25").Value = '  -2.85%  '
$ws.Range("D26").Value = '''163.40'
$ws.Range("E26").Value = '  -0.12%  '
$ws.Range("D27").Value = '''18.47'
$ws.Range("E27").Value = '  +1.24%  '
$ws.Range("D28").Value = '''1.913'
$ws.Range("E28").Value = '  -0.91%  '
$ws.Range("D29").Value = '''1.428'
$ws.Range("E29").Value = '  -2.35%  '
$ws.Range("D30").Value = '''0.09193'
$ws.Range("E30").Value = '  +0.44%  '
$ws.Range("D31").Value = '''4.233'
$ws.Range("E31").Value = '  -3.19%  '
$ws.Range("D32").Value = '''3.950'
$ws.Range("E32").Value = '  -2.39%  '
$ws.Range("D33").Value = '''0.05011'
$ws.Range("E33").Value = '  -3.84%  '
$ws.Range("D34").Value = '''0.7388'
$ws.Range("E34").Value = '  -0.25%  '
$ws.Range("E35").Value = '  +3.24%  '
$ws.Range("D36").Value = '''2.686'
$ws.Range("E36").Value = '  -1.13%  '
$ws.Range("D37").Value = '''0.01815'
$ws.Range("E37").Value = '  -0.01%  '
$ws.Range("D38").Value = '''2.602'
$ws.Range("E38").Value = '  -2.72%  '
$ws.Range("D39").Value = '''0.9039'
$ws.Range("E39").Value = '  -1.75%  '
$ws.Range("D40").Value = '''2.052'
$ws.Range("E40").Value = '  +0.96%  '
$ws.Range("D41").Value = '''5.937'
$ws.Range("E41").Value = '  -0.05%  '
$ws.Range("E42").Value = '  +0.48%  '
$ws.Range("D43").Value = '''0.4246'
$ws.Range("E43").Value = '  -3.51%  '
$ws.Range("D44").Value = '''0.9986'
$ws.Range("E44").Value = '  +0.56%  '
$ws.Range("D45").Value = '''7.392'
$ws.Range("E45").Value = '  -2.35%  '
$ws.Range("D46").Value = '''0.1311'
$ws.Range("E46").Value = '  -4.46%  '
$ws.Range("D47").Value = '''1.563'
$ws.Range("E47").Value = '  +11.18%  '
$ws.Range("D48").Value = '''63.96'
$ws.Range("E48").Value = '  -6.45%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '''8.774'
$ws.Range("E49").Value = '  -2.27%  '
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").Value = '''34.25'
$ws.Range("E50").Value = '  -2.05%  '
$ws.Range("D51").Value = '''0.05662'
$ws.Range("E51").Value = '  -2.72%  '
